# edit.ps1 - apply the README.docx changes described in the commit
# "Adding tweets for 14th and testing rbind":
#
#   1. "Jason " / "Bruiners" (split across two runs with a spellcheck
#      proofErr wrapper) becomes a single run "Jason Bruiners".
#   2. A new paragraph "Describe What each rmd contains" is inserted
#      right before the "Search Terms used:" paragraph.
#   3. The q1 search-term paragraph (several runs broken up by
#      proofErr spell-check markers) becomes one single run.
#   4. The q2 search-term paragraph (main text run + a trailing
#      single-space run) becomes one single run, trailing space
#      included, with xml:space="preserve".
#
# Because items 1, 3 and 4 do not actually change the *text* of the
# paragraph (only how it is chopped up into runs / proofErr markers),
# a plain Find-and-Replace is a no-op as far as the engine is
# concerned and leaves stray <w:proofErr/> markers behind. So instead,
# for each such paragraph we:
#   - insert a brand new, plain paragraph immediately before it,
#   - type the desired (clean) text into that new paragraph,
#   - delete the whole original paragraph (text + mark).
# That guarantees the final paragraph is a single plain run with no
# leftover proofErr markers or stale formatting.

$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1) "Jason " + proofErr + "Bruiners" + proofErr  ->  "Jason Bruiners"
# ---------------------------------------------------------------
$jasonOld = $d.Paragraphs(1)
$jasonOld.Range.InsertParagraphBefore()
$jasonNew = $d.Paragraphs(1)
$jasonNew.Range.InsertBefore("Jason Bruiners")
$jasonOld = $d.Paragraphs(2)
$jasonOld.Range.Delete()

# ---------------------------------------------------------------
# 2) Insert the new "Describe What each rmd contains" paragraph
#    right before "Search Terms used:"
# ---------------------------------------------------------------
$searchTerms = $d.Content
$searchTerms.Find.Execute("Search Terms used:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$searchTerms.InsertParagraphBefore()

# The empty paragraph we just inserted sits immediately before
# "Search Terms used:" - locate it via the document's paragraph list
# and fill it in with the new text.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs($i)
    if ($cand.Range.Text -eq "Search Terms used:" + [char]13) {
        $descParaIndex = $i - 1
        break
    }
}
$descPara = $d.Paragraphs($descParaIndex)
$descPara.Range.InsertBefore("Describe What each rmd contains")

# ---------------------------------------------------------------
# 3) Merge the q1 paragraph's runs/proofErr markers into one run
# ---------------------------------------------------------------
$q1Text = 'q1 <- "Covid OR Corona OR Covid-19 OR Covid19 OR Coronavirus OR SARS-CoV-2 OR SARS-CoV"'
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs($i)
    if ($cand.Range.Text -eq $q1Text + [char]13) {
        $q1Index = $i
        break
    }
}
$q1Old = $d.Paragraphs($q1Index)
$q1Old.Range.InsertParagraphBefore()
$q1New = $d.Paragraphs($q1Index)
$q1New.Range.InsertBefore($q1Text)
$q1Old = $d.Paragraphs($q1Index + 1)
$q1Old.Range.Delete()

# ---------------------------------------------------------------
# 4) Merge the q2 paragraph's two runs (text + trailing space) into
#    one run, keeping the trailing space.
# ---------------------------------------------------------------
$q2Text = 'q2 <- "OR #Covid OR #Corona OR #Covid-19 OR #Covid19 OR #Coronavirus OR #SARS-CoV-2 OR #SARS-CoV" '
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs($i)
    if ($cand.Range.Text -eq $q2Text + [char]13) {
        $q2Index = $i
        break
    }
}
$q2Old = $d.Paragraphs($q2Index)
$q2Old.Range.InsertParagraphBefore()
$q2New = $d.Paragraphs($q2Index)
$q2New.Range.InsertBefore($q2Text)
$q2Old = $d.Paragraphs($q2Index + 1)
$q2Old.Range.Delete()
